$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 1: "Nghiem trong (Critical)" severity bullet
#   3 Loi (30%) -> 2 Loi (20%)
#   remove "loi chen ma doc (XSS), " from the explanation sentence
# -----------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Nghiêm trọng (Critical) - 3 Lỗi (30%):")
$r.Text = "Nghiêm trọng (Critical) - 2 Lỗi (20%):"

$r = $d.Content
$r.Find.Execute(" Các lỗi đe dọa trực tiếp đến an toàn dữ liệu và doanh thu, bao gồm lỗi chèn mã độc (XSS), cho phép nhập số lượng âm vào giỏ hàng và thanh toán thẻ tín dụng không cần xác thực.")
$r.Text = " Các lỗi đe dọa trực tiếp đến an toàn dữ liệu và doanh thu, bao gồm lỗi cho phép nhập số lượng âm vào giỏ hàng và thanh toán thẻ tín dụng không cần xác thực."

# -----------------------------------------------------------------------
# Edit 2: "Lon (Major)" severity bullet
#   4 Loi (40%) -> 5 Loi (40%)
#   append ", loi hien thi ket qua tim kiem" before the closing parenthesis
# -----------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Lớn (Major) - 4 Lỗi (40%):")
$r.Text = "Lớn (Major) - 5 Lỗi (40%):"

$r = $d.Content
$r.Find.Execute(" Các lỗi về logic nghiệp vụ và thiếu hụt tính năng (Không có nút Quên mật khẩu, không có trang Lịch sử đơn hàng, lỗi cho phép nhập mật khẩu quá ngắn).")
$r.Text = " Các lỗi về logic nghiệp vụ và thiếu hụt tính năng (Không có nút Quên mật khẩu, không có trang Lịch sử đơn hàng, lỗi cho phép nhập mật khẩu quá ngắn, lỗi hiển thị kết quả tìm kiếm)."

# -----------------------------------------------------------------------
# Edit 3: "Top 3 loi Nghiem trong (Critical)" -> "Top 3 loi noi com nhat (Outstanding Issues)"
# -----------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Top 3 lỗi Nghiêm trọng (Critical)")
$r.Text = "Top 3 lỗi nổi cộm nhất (Outstanding Issues)"

# -----------------------------------------------------------------------
# Edit 4: BUG_PROD_001 - replace the invalid XSS bug description with the
# correct UX / missing-search-keyword issue description.
# -----------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("[BUG_PROD_001] Lỗ hổng bảo mật XSS trên thanh Tìm kiếm:")
$r.Text = "[BUG_PROD_001] Lỗi Trải nghiệm người dùng (UX) tại trang kết quả tìm kiếm: "

$r = $d.Content
$r.Find.Execute(" Hệ thống không mã hóa (Encode) các thẻ HTML/Javascript được nhập vào thanh tìm kiếm, tạo ra lỗ hổng để tin tặc đánh cắp phiên đăng nhập của người dùng.")
$r.Text = "Sau khi thực hiện truy vấn, hệ thống chỉ trả về một tiêu đề tĩnh là ""SEARCHED PRODUCTS"" mà không hiển thị lại từ khóa mà người dùng vừa nhập. Điều này vi phạm nguyên tắc cơ bản về tính khả dụng (Usability) trong thiết kế, gây bối rối cho khách hàng khi họ cần đối chiếu xem mình có gõ sai chính tả hay không."
$r.Bold = 0
